# "systeme de notification intelligente"
#
# 1) The footer/date placeholder on the slide master + every slide layout
#    shows a fixed "update automatically" date field (type="datetimeFigureOut").
#    Its cached text moves from 17/03/2025 -> 21/03/2025.
# 2) The slide 1 title is retouched: "Serviance" -> "Serviellance".
#    The title placeholder's body also picks up a normAutofit shrink
#    (fontScale) because the longer run no longer fits on one line.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- Slide master date placeholder ---
Set-DatePlaceholderText -shapes $p.SlideMaster.Shapes -newText "21/03/2025"

# --- Every slide layout's date placeholder ---
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Set-DatePlaceholderText -shapes $layout.Shapes -newText "21/03/2025"
}

# --- Slide 1 title text fix: "Serviance" -> "Serviellance" ---
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Title
$titleRange = $title.TextFrame.TextRange
$fullText = $titleRange.Text
$oldWord = "Serviance"
$newWord = "Serviellance"
$startPos = $fullText.IndexOf($oldWord) + 1

if ($startPos -gt 0) {
    $wordRange = $titleRange.Characters($startPos, $oldWord.Length)
    $wordRange.Text = $newWord
}

# The title text box uses normAutofit; once the replacement text is a
# little longer it no longer fits on a single line at 100%, so PowerPoint
# shrinks it (fontScale=90%) to keep it inside the placeholder.
try { $title.TextFrame.AutoSize = 2 } catch {}
try { $title.TextFrame.AutofitFontScale = 90000 } catch {}
try { $title.TextFrame2.FontScale = 90000 } catch {}
